$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows whose coin identity is unchanged.
# A leading apostrophe is used on column D assignments so Excel stores the dotted
# price strings as text (matching the source data) instead of auto-converting to numbers.
$ws.Range("D2").Value = "'27.375.40"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "'1.860.70"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "'314.55"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "'0.4644"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.3718"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").Value = "'0.07345"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").Value = "'0.8847"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("D11").Value = "'0.07901"
$ws.Range("E11").Value = "  +3.84%  "
$ws.Range("D12").Value = "'19.89"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "'1.877.13"
$ws.Range("E13").Value = "  +4.01%  "
$ws.Range("D14").Value = "'5.404"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "'6.581"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "'92.18"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'0.000008875"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "'14.84"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "'27.425.07"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").Value = "'5.142"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'10.54"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'2.145.38"
$ws.Range("E24").Value = "  +6.67%  "
$ws.Range("D25").Value = "'1.901"
$ws.Range("E25").Value = "  +3.07%  "
$ws.Range("D26").Value = "'152.93"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "'18.42"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").Value = "'2.069"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "'5.126"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").Value = "'116.39"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").Value = "'0.08898"
$ws.Range("E31").Value = "  +0.68%  "

# Rows 32/33: HuobiToken and ImmutableX swap positions (row 32 becomes HuobiToken, row 33 becomes ImmutableX)
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "'3.026"
$ws.Range("E32").Value = "  +2.46%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7539"
$ws.Range("E33").Value = "  +4.76%  "

$ws.Range("D34").Value = "'1.161"
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("D35").Value = "'4.493"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "'2.642"
$ws.Range("E36").Value = "  +9.69%  "
$ws.Range("D37").Value = "'0.01965"
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("D38").Value = "'1.080"
$ws.Range("E38").Value = "  +0.48%  "

# Rows 39/40: Hedera and MXToken swap positions (row 39 becomes Hedera, row 40 becomes MXToken)
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05252"
$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.976"
$ws.Range("E40").Value = "  +1.61%  "

$ws.Range("D41").Value = "'7.122"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").Value = "'0.5162"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").Value = "'8.346"
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("D45").Value = "'0.4841"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "'10.31"
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "'103.69"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").Value = "'1.651"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").Value = "'0.06245"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "'65.85"
$ws.Range("E51").Value = "  +2.63%  "
